$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"224.7833046666667"
$ws.Range("H2").Value = [double]"674.349914"
$ws.Range("I2").Value = [double]"0.3882379172278888"
$ws.Range("J2").Value = [double]"0.3882379172278889"
$ws.Range("M2").Value = [double]"3.425446666666666"
$ws.Range("N2").Value = [double]"10.27634"
$ws.Range("O2").Value = [double]"0.6657953389778073"
$ws.Range("P2").Value = [double]"0.6657953389778073"
$ws.Range("Q2").Value = [double]"769.983221692751"
$ws.Range("R2").Value = [double]"6929.84899523476"
$ws.Range("S2").Value = [double]"0.2584869957047801"
$ws.Range("T2").Value = [double]"0.2584869957047802"
$ws.Range("G3").Value = [double]"224.7833046666667"
$ws.Range("H3").Value = [double]"674.349914"
$ws.Range("I3").Value = [double]"0.3882379172278888"
$ws.Range("J3").Value = [double]"0.3882379172278889"
$ws.Range("O3").Value = [double]"0.2094791321596951"
$ws.Range("P3").Value = [double]"0.2094791321596952"
$ws.Range("Q3").Value = [double]"242.2597570378907"
$ws.Range("R3").Value = [double]"2180.337813341016"
$ws.Range("S3").Value = [double]"0.08132774197238569"
$ws.Range("T3").Value = [double]"0.08132774197238572"
$ws.Range("G4").Value = [double]"224.7833046666667"
$ws.Range("H4").Value = [double]"674.349914"
$ws.Range("I4").Value = [double]"0.3882379172278888"
$ws.Range("J4").Value = [double]"0.3882379172278889"
$ws.Range("M4").Value = [double]"0.62317"
$ws.Range("N4").Value = [double]"1.86951"
$ws.Range("O4").Value = [double]"0.1211239647746572"
$ws.Range("P4").Value = [double]"0.1211239647746572"
$ws.Range("Q4").Value = [double]"140.0782119691267"
$ws.Range("R4").Value = [double]"1260.70390772214"
$ws.Range("S4").Value = [double]"0.04702491581049708"
$ws.Range("T4").Value = [double]"0.04702491581049708"
$ws.Range("G5").Value = [double]"224.7833046666667"
$ws.Range("H5").Value = [double]"674.349914"
$ws.Range("I5").Value = [double]"0.3882379172278888"
$ws.Range("J5").Value = [double]"0.3882379172278889"
$ws.Range("M5").Value = [double]"0.01852966666666667"
$ws.Range("N5").Value = [double]"0.055589"
$ws.Range("O5").Value = [double]"0.003601564087840353"
$ws.Range("P5").Value = [double]"0.003601564087840353"
$ws.Range("Q5").Value = [double]"4.165159707705111"
$ws.Range("R5").Value = [double]"37.486437369346"
$ws.Range("S5").Value = [double]"0.0013982637402259"
$ws.Range("T5").Value = [double]"0.0013982637402259"
$ws.Range("I6").Value = [double]"0.4251955538547045"
$ws.Range("J6").Value = [double]"0.4251955538547046"
$ws.Range("M6").Value = [double]"3.425446666666666"
$ws.Range("N6").Value = [double]"10.27634"
$ws.Range("O6").Value = [double]"0.6657953389778073"
$ws.Range("P6").Value = [double]"0.6657953389778073"
$ws.Range("Q6").Value = [double]"843.2804419108421"
$ws.Range("R6").Value = [double]"7589.523977197579"
$ws.Range("S6").Value = [double]"0.2830932179105495"
$ws.Range("T6").Value = [double]"0.2830932179105495"
$ws.Range("I7").Value = [double]"0.4251955538547045"
$ws.Range("J7").Value = [double]"0.4251955538547046"
$ws.Range("O7").Value = [double]"0.2094791321596951"
$ws.Range("P7").Value = [double]"0.2094791321596952"
$ws.Range("S7").Value = [double]"0.08906959561964442"
$ws.Range("T7").Value = [double]"0.08906959561964445"
$ws.Range("I8").Value = [double]"0.4251955538547045"
$ws.Range("J8").Value = [double]"0.4251955538547046"
$ws.Range("M8").Value = [double]"0.62317"
$ws.Range("N8").Value = [double]"1.86951"
$ws.Range("O8").Value = [double]"0.1211239647746572"
$ws.Range("P8").Value = [double]"0.1211239647746572"
$ws.Range("Q8").Value = [double]"153.4127149312633"
$ws.Range("R8").Value = [double]"1380.71443438137"
$ws.Range("S8").Value = [double]"0.05150137128743808"
$ws.Range("T8").Value = [double]"0.05150137128743809"
$ws.Range("I9").Value = [double]"0.4251955538547045"
$ws.Range("J9").Value = [double]"0.4251955538547046"
$ws.Range("M9").Value = [double]"0.01852966666666667"
$ws.Range("N9").Value = [double]"0.055589"
$ws.Range("O9").Value = [double]"0.003601564087840353"
$ws.Range("P9").Value = [double]"0.003601564087840353"
$ws.Range("Q9").Value = [double]"4.561654877649223"
$ws.Range("R9").Value = [double]"41.054893898843"
$ws.Range("S9").Value = [double]"0.001531369037072493"
$ws.Range("T9").Value = [double]"0.001531369037072493"
$ws.Range("G10").Value = [double]"107.695137"
$ws.Range("H10").Value = [double]"323.085411"
$ws.Range("I10").Value = [double]"0.186007300437435"
$ws.Range("J10").Value = [double]"0.186007300437435"
$ws.Range("M10").Value = [double]"3.425446666666666"
$ws.Range("N10").Value = [double]"10.27634"
$ws.Range("O10").Value = [double]"0.6657953389778073"
$ws.Range("P10").Value = [double]"0.6657953389778073"
$ws.Range("Q10").Value = [double]"368.90394805286"
$ws.Range("R10").Value = [double]"3320.13553247574"
$ws.Range("S10").Value = [double]"0.1238427936470889"
$ws.Range("T10").Value = [double]"0.1238427936470889"
$ws.Range("G11").Value = [double]"107.695137"
$ws.Range("H11").Value = [double]"323.085411"
$ws.Range("I11").Value = [double]"0.186007300437435"
$ws.Range("J11").Value = [double]"0.186007300437435"
$ws.Range("O11").Value = [double]"0.2094791321596951"
$ws.Range("P11").Value = [double]"0.2094791321596952"
$ws.Range("Q11").Value = [double]"116.068218511476"
$ws.Range("R11").Value = [double]"1044.613966603284"
$ws.Range("S11").Value = [double]"0.03896464787100155"
$ws.Range("T11").Value = [double]"0.03896464787100157"
$ws.Range("G12").Value = [double]"107.695137"
$ws.Range("H12").Value = [double]"323.085411"
$ws.Range("I12").Value = [double]"0.186007300437435"
$ws.Range("J12").Value = [double]"0.186007300437435"
$ws.Range("M12").Value = [double]"0.62317"
$ws.Range("N12").Value = [double]"1.86951"
$ws.Range("O12").Value = [double]"0.1211239647746572"
$ws.Range("P12").Value = [double]"0.1211239647746572"
$ws.Range("Q12").Value = [double]"67.11237852429001"
$ws.Range("R12").Value = [double]"604.01140671861"
$ws.Range("S12").Value = [double]"0.02252994170601295"
$ws.Range("T12").Value = [double]"0.02252994170601295"
$ws.Range("G13").Value = [double]"107.695137"
$ws.Range("H13").Value = [double]"323.085411"
$ws.Range("I13").Value = [double]"0.186007300437435"
$ws.Range("J13").Value = [double]"0.186007300437435"
$ws.Range("M13").Value = [double]"0.01852966666666667"
$ws.Range("N13").Value = [double]"0.055589"
$ws.Range("O13").Value = [double]"0.003601564087840353"
$ws.Range("P13").Value = [double]"0.003601564087840353"
$ws.Range("Q13").Value = [double]"1.995554990231"
$ws.Range("R13").Value = [double]"17.959994912079"
$ws.Range("S13").Value = [double]"0.0006699172133315969"
$ws.Range("T13").Value = [double]"0.0006699172133315971"
$ws.Range("G14").Value = [double]"0.323784"
$ws.Range("H14").Value = [double]"0.971352"
$ws.Range("I14").Value = [double]"0.0005592284799715185"
$ws.Range("J14").Value = [double]"0.0005592284799715186"
$ws.Range("M14").Value = [double]"3.425446666666666"
$ws.Range("N14").Value = [double]"10.27634"
$ws.Range("O14").Value = [double]"0.6657953389778073"
$ws.Range("P14").Value = [double]"0.6657953389778073"
$ws.Range("Q14").Value = [double]"1.10910482352"
$ws.Range("R14").Value = [double]"9.98194341168"
$ws.Range("S14").Value = [double]"0.0003723317153886811"
$ws.Range("T14").Value = [double]"0.0003723317153886812"
$ws.Range("G15").Value = [double]"0.323784"
$ws.Range("H15").Value = [double]"0.971352"
$ws.Range("I15").Value = [double]"0.0005592284799715185"
$ws.Range("J15").Value = [double]"0.0005592284799715186"
$ws.Range("O15").Value = [double]"0.2094791321596951"
$ws.Range("P15").Value = [double]"0.2094791321596952"
$ws.Range("Q15").Value = [double]"0.348957558432"
$ws.Range("R15").Value = [double]"3.140618025888"
$ws.Range("S15").Value = [double]"0.0001171466966634192"
$ws.Range("T15").Value = [double]"0.0001171466966634192"
$ws.Range("G16").Value = [double]"0.323784"
$ws.Range("H16").Value = [double]"0.971352"
$ws.Range("I16").Value = [double]"0.0005592284799715185"
$ws.Range("J16").Value = [double]"0.0005592284799715186"
$ws.Range("M16").Value = [double]"0.62317"
$ws.Range("N16").Value = [double]"1.86951"
$ws.Range("O16").Value = [double]"0.1211239647746572"
$ws.Range("P16").Value = [double]"0.1211239647746572"
$ws.Range("Q16").Value = [double]"0.20177247528"
$ws.Range("R16").Value = [double]"1.81595227752"
$ws.Range("S16").Value = [double]"6.773597070905529E-05"
$ws.Range("T16").Value = [double]"6.773597070905531E-05"
$ws.Range("G17").Value = [double]"0.323784"
$ws.Range("H17").Value = [double]"0.971352"
$ws.Range("I17").Value = [double]"0.0005592284799715185"
$ws.Range("J17").Value = [double]"0.0005592284799715186"
$ws.Range("M17").Value = [double]"0.01852966666666667"
$ws.Range("N17").Value = [double]"0.055589"
$ws.Range("O17").Value = [double]"0.003601564087840353"
$ws.Range("P17").Value = [double]"0.003601564087840353"
$ws.Range("Q17").Value = [double]"0.005999609592000001"
$ws.Range("R17").Value = [double]"0.053996486328"
$ws.Range("S17").Value = [double]"2.014097210362969E-06"
$ws.Range("T17").Value = [double]"2.014097210362969E-06"
